$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D7 ("宝箱" row): effect text replaced with a new "gain loot" effect.
$ws.Range("D7").Value = "获得遭遇牌堆第1张战利品牌。使用1张《敏捷》发动本牌时，可以再获得1张战利品牌。"

# D8 ("冒险者尸体" row): "遗物牌堆顶的1张" -> "遗物牌堆第1张".
$ws.Range("D8").Value = "消耗3时间，将主牌堆第1张怪物牌放在房间区任意空槽位，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆第1张遗物牌。"

# View/scroll position changed so row 6 is at the top of the visible area,
# and the active selection moved from D12 to D10.
$excel.Goto($ws.Range("A6"), $true)
$ws.Range("D10").Select()
